$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the "Created" date
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2020-04-08", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2020-05-04", 2)

# ---------------------------------------------------------------------------
# 2. Create the new custom styles referenced by the new content
# ---------------------------------------------------------------------------
function New-HiddenCharStyle($styleId) {
    $s = $d.Styles.Add($styleId, 2)
    $s.BaseStyle = $d.Styles("BodyTextChar")
    $s.Visibility = $true
    return $s
}

function New-HiddenParaStyle($styleId) {
    $s = $d.Styles.Add($styleId, 1)
    $s.BaseStyle = $d.Styles("BodyText")
    $s.QuickStyle = $true
    $s.Visibility = $true
    return $s
}

foreach ($i in 9..15) {
    New-HiddenCharStyle("redoc-inlinecode-$i") | Out-Null
}
New-HiddenParaStyle("redoc-codechunk-3") | Out-Null
New-HiddenParaStyle("redoc-codechunk-4") | Out-Null

# ---------------------------------------------------------------------------
# Helper: append a run of text at $pos, optionally styled + bookmarked.
# Returns the end position after the inserted text.
# ---------------------------------------------------------------------------
function Add-Run($pos, $text, $styleName, $bookmarkName) {
    $rng = $d.Range($pos, $pos)
    $rng.InsertAfter($text)
    $endPos = $pos + $text.Length
    if ($styleName) {
        $styled = $d.Range($pos, $endPos)
        $styled.Style = $d.Styles($styleName)
        if ($bookmarkName) {
            $d.Bookmarks.Add($bookmarkName, $styled)
        }
    }
    return $endPos
}

# Helper: insert a brand-new empty paragraph at the end of the document,
# return its start position (where its (still empty) text begins).
function Add-EmptyParagraphAtEnd() {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    return $r.Start
}

# Helper: create a hidden "code chunk" paragraph (empty, with the given
# paragraph style) at the end of the document, wrapped by a bookmark with
# the same name as the style.
function Add-CodeChunkParagraph($styleName) {
    $pos = Add-EmptyParagraphAtEnd
    $para = $d.Paragraphs.Last
    $para.Range.Style = $d.Styles($styleName)

    # A bookmark anchored on a truly empty range is unreliable, so we
    # insert a placeholder character, bookmark it, then remove it again.
    $tmp = $d.Range($pos, $pos)
    $tmp.InsertAfter(" ")
    $tmpRng = $d.Range($pos, $pos + 1)
    $d.Bookmarks.Add($styleName, $tmpRng)
    $clearRng = $d.Range($pos, $pos + 1)
    $clearRng.Text = ""
}

# ---------------------------------------------------------------------------
# 3. Replace the content of the last ("In a regional perspective...")
#    paragraph with the new, more detailed text.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$para9 = $d.Paragraphs.Item($lastParaIndex)
$para9.Range.Text = ""
$pos = $para9.Range.Start

$pos = Add-Run $pos "In a regional perspective, the bottom 40 of most of the countries of" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "East Asia and Pacific, Europe and Central Asia, and Latin America and the Caribbean" "redoc-inlinecode-9" "redoc-inlinecode-9"
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "increased their income over the average; meanwhile, in" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "South Asia, Middle East and North Africa, Sub-saharan Africa, and Other High Income countries" "redoc-inlinecode-10" "redoc-inlinecode-10"
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "their income improved less than the national average, being" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "100%, 75%, 66.67% and 54.55%" "redoc-inlinecode-11" "redoc-inlinecode-11"
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "of their countries exhibiting such behaiviour, respectively." $null $null

# ---------------------------------------------------------------------------
# 4. Add the new "An overview by income groups..." paragraph.
# ---------------------------------------------------------------------------
$pos = Add-EmptyParagraphAtEnd
$d.Paragraphs.Last.Range.Style = $d.Styles("BodyText")

$pos = Add-Run $pos "An overview by income groups is largly meanigfull. The" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "57.14% and 79.17%" "redoc-inlinecode-12" "redoc-inlinecode-12"
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "of" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "Upper middle income, and High income" "redoc-inlinecode-13" "redoc-inlinecode-13"
$pos = Add-Run $pos "’s bottom 40% increased their income more rapidly than de average, accordingly; Whereas, for the" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "Low income, and Lower middle income" "redoc-inlinecode-14" "redoc-inlinecode-14"
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "countries the bottom earners experienced a lower expansion of their wellfare compared to the average, with only" $null $null
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "33.33% and 45.45%" "redoc-inlinecode-15" "redoc-inlinecode-15"
$pos = Add-Run $pos " " $null $null
$pos = Add-Run $pos "of the countries, in the respective income group, with the bottom’s growth out-performing the average." $null $null

# ---------------------------------------------------------------------------
# 5. Hidden code-chunk paragraph #3
# ---------------------------------------------------------------------------
Add-CodeChunkParagraph "redoc-codechunk-3"

# ---------------------------------------------------------------------------
# 6. Add the "It can be valuable..." paragraph.
# ---------------------------------------------------------------------------
$pos = Add-EmptyParagraphAtEnd
$d.Paragraphs.Last.Range.Style = $d.Styles("BodyText")
$pos = Add-Run $pos "It can be valuable to check the how the countries have performed across time. Comparing the current growth rates (2012-2017) with the ones between 2008-2013, we found broad changes in the countries behaivior. For both growth rate windows, the following figure shows how the contries rank form the one with the fastest growth of the bottom compared to the average, to the one with the slowest." $null $null

# ---------------------------------------------------------------------------
# 7. Hidden code-chunk paragraph #4
# ---------------------------------------------------------------------------
Add-CodeChunkParagraph "redoc-codechunk-4"

Write-Output "Edit complete"
